$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D and E) for ownTeam/oppTeam, shifting existing
# batsman..sr columns from D..I to F..K
$ws.Range("D1:E1").EntireColumn.Insert()

# Update header row for the two newly inserted columns
$ws.Cells.Item(1, 4).Value2 = "ownTeam"
$ws.Cells.Item(1, 5).Value2 = "oppTeam"

# Full data table (rows 2-18), columns A-K:
# venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr
$data = @(
    @(' Dubai (DSC)', ' October 27 2020', 'Sunrisers won by 88 runs', 'Delhi Capitals', 'Sunrisers Hyderabad', 'Shreyas Iyer (c)', '''7', '''12', '''0', '''0', '''58.33'),
    @(' Dubai (DSC)', ' November 05 2020', 'Mumbai won by 57 runs', 'Delhi Capitals', 'Mumbai Indians', 'Shreyas Iyer (c)', '''12', '''8', '''3', '''0', '''150.00'),
    @(' Dubai (DSC)', ' November 10 2020', 'Mumbai won by 5 wickets (with 8 balls remaining)', 'Delhi Capitals', 'Mumbai Indians', 'Shreyas Iyer (c)', '''65', '''50', '''6', '''2', '''130.00'),
    @(' Dubai (DSC)', ' October 14 2020', 'Capitals won by 13 runs', 'Delhi Capitals', 'Rajasthan Royals', 'Shreyas Iyer (c)', '''53', '''43', '''3', '''2', '''123.25'),
    @(' Dubai (DSC)', ' October 31 2020', 'Mumbai won by 9 wickets (with 34 balls remaining)', 'Delhi Capitals', 'Mumbai Indians', 'Shreyas Iyer (c)', '''25', '''29', '''1', '''1', '''86.20'),
    @(' Dubai (DSC)', ' September 25 2020', 'Capitals won by 44 runs', 'Delhi Capitals', 'Chennai Super Kings', 'Shreyas Iyer (c)', '''26', '''22', '''1', '''0', '''118.18'),
    @(' Sharjah', ' October 17 2020', 'Capitals won by 5 wickets (with 1 ball remaining)', 'Delhi Capitals', 'Chennai Super Kings', 'Shreyas Iyer (c)', '''23', '''23', '''1', '''1', '''100.00'),
    @(' Sharjah', ' October 09 2020', 'Capitals won by 46 runs', 'Delhi Capitals', 'Rajasthan Royals', 'Shreyas Iyer (c)', '''22', '''18', '''4', '''0', '''122.22'),
    @(' Dubai (DSC)', ' October 05 2020', 'Capitals won by 59 runs', 'Delhi Capitals', 'Royal Challengers Bangalore', 'Shreyas Iyer (c)', '''11', '''13', '''1', '''0', '''84.61'),
    @(' Dubai (DSC)', ' October 20 2020', 'Kings XI won by 5 wickets (with 6 balls remaining)', 'Delhi Capitals', 'Kings XI Punjab', 'Shreyas Iyer (c)', '''14', '''12', '''0', '''1', '''116.66'),
    @(' Dubai (DSC)', ' September 20 2020', 'Match tied (Capitals won the one-over eliminator)', 'Delhi Capitals', 'Kings XI Punjab', 'Shreyas Iyer (c)', '''39', '''32', '''0', '''3', '''121.87'),
    @(' Sharjah', ' October 03 2020', 'Capitals won by 18 runs', 'Delhi Capitals', 'Kolkata Knight Riders', 'Shreyas Iyer (c)', '''88', '''38', '''7', '''6', '''231.57'),
    @(' Abu Dhabi', ' November 08 2020', 'Capitals won by 17 runs', 'Delhi Capitals', 'Sunrisers Hyderabad', 'Shreyas Iyer (c)', '''21', '''20', '''1', '''0', '''105.00'),
    @(' Abu Dhabi', ' November 02 2020', 'Capitals won by 6 wickets (with 6 balls remaining)', 'Delhi Capitals', 'Royal Challengers Bangalore', 'Shreyas Iyer (c)', '''7', '''9', '''0', '''0', '''77.77'),
    @(' Abu Dhabi', ' October 24 2020', 'KKR won by 59 runs', 'Delhi Capitals', 'Kolkata Knight Riders', 'Shreyas Iyer (c)', '''47', '''38', '''5', '''0', '''123.68'),
    @(' Abu Dhabi', ' October 11 2020', 'Mumbai won by 5 wickets (with 2 balls remaining)', 'Delhi Capitals', 'Mumbai Indians', 'Shreyas Iyer (c)', '''42', '''33', '''5', '''0', '''127.27'),
    @(' Abu Dhabi', ' September 29 2020', 'Sunrisers won by 15 runs', 'Delhi Capitals', 'Sunrisers Hyderabad', 'Shreyas Iyer (c)', '''17', '''21', '''2', '''0', '''80.95')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $rowNum = $i + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value2 = $row[$c]
    }
}

Write-Output "done"
